$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.680.81'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '2.244.63'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.633'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.77%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0933'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.893'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.57%  '
$ws.Range('D16').Value = '2.580.92'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '2.249.38'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '42.738.04'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.60%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000108'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +19.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '232.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.32%  '
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.03%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '175.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0911'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +21.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0374'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('E41').Value = '  +4.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.236'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.26%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.651'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.39%  '
